$wb = $excel.ActiveWorkbook

# ---- Sheet "Summary" ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5071174377224199
$ws1.Range("C2").Value = 0.08080808080808081
$ws1.Range("D2").Value = 0.8571428571428571
$ws1.Range("E2").Value = 0.1476923076923077
$ws1.Range("F2").Value = 0.293398533007335
$ws1.Range("G2").Value = 0.6258776328986961
$ws1.Range("H2").Value = 0.6858614232209738
$ws1.Range("I2").Value = 24
$ws1.Range("J2").Value = 273
$ws1.Range("K2").Value = 261
$ws1.Range("L2").Value = 4

# ---- Sheet "Classification Report" ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 - label "0"
$ws2.Range("B2").Value = 0.9849056603773585
$ws2.Range("C2").Value = 0.4887640449438202
$ws2.Range("D2").Value = 0.6533166458072591

# Row 3 - label "1"
$ws2.Range("B3").Value = 0.08080808080808081
$ws2.Range("C3").Value = 0.8571428571428571
$ws2.Range("D3").Value = 0.1476923076923077

# Row 4 - label "accuracy"
$ws2.Range("B4").Value = 0.5071174377224199
$ws2.Range("C4").Value = 0.5071174377224199
$ws2.Range("D4").Value = 0.5071174377224199
$ws2.Range("E4").Value = 0.5071174377224199

# Row 5 - label "macro avg"
$ws2.Range("B5").Value = 0.5328568705927197
$ws2.Range("C5").Value = 0.6729534510433386
$ws2.Range("D5").Value = 0.4005044767497834

# Row 6 - label "weighted avg"
$ws2.Range("B6").Value = 0.939861652854334
$ws2.Range("C6").Value = 0.5071174377224199
$ws2.Range("D6").Value = 0.6281253976449483

# ---- Sheet "Confusion Matrix" ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - "Actual 0"
$ws3.Range("B2").Value = 261
$ws3.Range("C2").Value = 273

# Row 3 - "Actual 1"
$ws3.Range("B3").Value = 4
$ws3.Range("C3").Value = 24
